# Add a new row of scanned-data records (row 57) to each of the four
# worksheets in the workbook. The new row follows the same layout/style
# as the existing rows (e.g. row 56): column A is a timestamp using the
# date-time style already applied to the other rows in column A, and the
# remaining columns hold the raw/decoded packet fields.

$wb = $excel.ActiveWorkbook

# New timestamp (column A) shared by the new record appended to each sheet.
$timestampValue = 45843.43393518519

$rowsData = @{
    "DE_LFT_#1" = @{
        A = $timestampValue
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x5C"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 348
        I = 14
    }
    "DE_LFT_#2" = @{
        A = $timestampValue
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x60"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 352
        I = 14
    }
    "DE_PLT_#1" = @{
        A = $timestampValue
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7E"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 126
        I = 7
    }
    "DE_PLT_#2" = @{
        A = $timestampValue
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7D"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 125
        I = 3
    }
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Each sheet currently has data through row 56; append the new record
    # as row 57 (and update the sheet's used-range dimension to match).
    $lastRow = $ws.UsedRange.Rows.Count
    $row = $lastRow + 1

    $data = $rowsData[$sheetName]

    # Column A: timestamp, copy number format from the row above it (A56)
    # so the new cell keeps the same date/time display style.
    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

    # Columns B-E: raw hex-byte strings stored as text.
    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E

    # Columns F-I: decoded numeric values.
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
}
